# Actualización 10 de Mayo
# Adds 3 rescatable students to the "Rescatables" worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A - control numbers
$ws.Range("A2").Value = 20330051920322
$ws.Range("A3").Value = 20330051920335
$ws.Range("A4").Value = 18330051920357

# Column B - Paterno (father's surname)
$ws.Range("B2").Value = "GARCIA"
$ws.Range("B3").Value = "VAZQUEZ"
$ws.Range("B4").Value = "BRAVO"

# Column C - Materno (mother's surname)
$ws.Range("C2").Value = "FLORES"
$ws.Range("C3").Value = "TZIZIHUA"
$ws.Range("C4").Value = "REYES"

# Column D - Nombres (given names)
$ws.Range("D2").Value = "MARCOS"
$ws.Range("D3").Value = "DORA LUZ"
$ws.Range("D4").Value = "PATRICIA MARLENE"

# Column E - Nombre_Largo (subject name)
$ws.Range("E2").Value = "ENSAMBLA E INSTALA CONTROLADORES Y DISPOSITIVOS PERIFÉRICOS"
$ws.Range("E3").Value = "ENSAMBLA E INSTALA CONTROLADORES Y DISPOSITIVOS PERIFÉRICOS"
$ws.Range("E4").Value = "DESARROLLA APLICACIONES MÓVILES PARA ANDROID"

# Column F - Grupo
$ws.Range("F2").Value = "2ASV"
$ws.Range("F3").Value = "2ASV"
$ws.Range("F4").Value = "6APM"

# Column G - Reprobadas
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
